$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 and Row 45 swap places (Mantle <-> Stellar), plus value updates throughout the table

# Row 44: becomes Stellar
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0978"
$ws.Range("E44").Value = "  -3.32%  "

# Row 45: becomes Mantle
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").Value = "  -2.93%  "

# Remaining price / volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.859.09"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.561.20"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.89"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.16"
$ws.Range("E6").Value = "  -7.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  -5.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.561.75"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.21"
$ws.Range("E10").Value = "  -7.28%  "
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  -5.08%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.005.01"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.855.50"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("E16").Value = "  -5.08%  "
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.558.82"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -5.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.85"
$ws.Range("E20").Value = "  -6.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.05"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.35"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").Value = "0.0₃0775"
$ws.Range("E28").Value = "  -8.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.86"
$ws.Range("E29").Value = "  -7.35%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.59"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.56"
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  -7.59%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -7.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.883"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -8.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.89"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  -9.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "286.75"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("E41").Value = "  -8.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  -5.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.61"
$ws.Range("E48").Value = "  -5.11%  "
$ws.Range("E50").Value = "  -8.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.912.72"
$ws.Range("E51").Value = "  -2.44%  "
